# Update the table style applied to every "Data Sources" style table in
# the deck: swap the old tableStyleId GUID for the new one.
#
# Old style: {C75BC050-3C19-41D3-958A-A709838A6AC1}
# New style: {86A1A7EE-9A7A-468C-9D75-9378E41A33B2}

$oldStyleId = "{C75BC050-3C19-41D3-958A-A709838A6AC1}"
$newStyleId = "{86A1A7EE-9A7A-468C-9D75-9378E41A33B2}"

$p = $ppt.ActivePresentation

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $shp = $s.Shapes.Item($i)
        if ($shp.HasTable) {
            $tbl = $shp.Table
            if ($tbl.Style -eq $oldStyleId) {
                $tbl.ApplyStyle($newStyleId)
            }
        }
    }
}
